$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) updates - force text format to preserve exact string
# representation (avoid Excel auto-converting numeric-looking strings to numbers)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.947.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.651.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0875"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.884.71"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.661.89"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.520"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.946.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "236.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.545.87"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.581"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.895"
$ws.Range("D38").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "66.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.792.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.938"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.75"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0987"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.64"
$ws.Range("D51").Style = "Normal"

# Volume(1h) percentage column (E) updates
$ws.Range("E2").Value = "  +2.13%  "
$ws.Range("E3").Value = "  +2.77%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  +1.31%  "
$ws.Range("E6").Value = "  +2.24%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +2.60%  "
$ws.Range("E9").Value = "  +1.70%  "
$ws.Range("E11").Value = "  +3.32%  "
$ws.Range("E12").Value = "  +2.84%  "
$ws.Range("E13").Value = "  +3.58%  "
$ws.Range("E14").Value = "  +2.24%  "
$ws.Range("E15").Value = "  +2.60%  "
$ws.Range("E16").Value = "  +2.88%  "
$ws.Range("E17").Value = "  +2.15%  "
$ws.Range("E18").Value = "  +2.33%  "
$ws.Range("E19").Value = "  +1.30%  "
$ws.Range("E20").Value = "  +1.73%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("E22").Value = "  +3.38%  "
$ws.Range("E23").Value = "  +3.81%  "
$ws.Range("E24").Value = "  +1.75%  "
$ws.Range("E25").Value = "  -1.21%  "
$ws.Range("E26").Value = "  +2.08%  "
$ws.Range("E27").Value = "  +0.89%  "
$ws.Range("E29").Value = "  +2.64%  "
$ws.Range("E30").Value = "  +0.50%  "
$ws.Range("E31").Value = "  +1.59%  "
$ws.Range("E32").Value = "  +3.20%  "
$ws.Range("E33").Value = "  +2.82%  "
$ws.Range("E34").Value = "  +4.75%  "
$ws.Range("E35").Value = "  +9.51%  "
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("E37").Value = "  +3.21%  "
$ws.Range("E39").Value = "  +2.69%  "
$ws.Range("E40").Value = "  +3.57%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("E42").Value = "  +8.81%  "
$ws.Range("E43").Value = "  +2.30%  "
$ws.Range("E44").Value = "  +2.73%  "
$ws.Range("E45").Value = "  +1.96%  "
$ws.Range("E46").Value = "  +0.94%  "
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("E48").Value = "  +2.06%  "
$ws.Range("E49").Value = "  +3.01%  "
$ws.Range("E50").Value = "  +0.93%  "
$ws.Range("E51").Value = "  +2.58%  "
